# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" colours, linked only from the
#                             notes master (ppt/notesMasters/notesMaster1.xml)
#   ppt/theme/theme2.xml  -> "Integral" colours, linked from the slide master
#                             (ppt/slideMasters/slideMaster1.xml), i.e. the
#                             theme that is actually applied to every slide.
#
# The authored change swaps the two themes' contents: the slide master's
# theme becomes the "Office Theme" colour palette (what used to live in
# theme1.xml), while the (unused-by-slides) notes-master theme becomes the
# "Integral" palette.  The part <-> relationship wiring itself is untouched;
# only the 12 scheme colours (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) that
# end up in the slide master's theme actually change in a way that affects
# the rendered deck.
#
# Apply this through the live PowerPoint colour-scheme object model so the
# host rewrites the underlying theme XML for us: walk the 12 theme colours
# via a slide's ThemeColorScheme and set each RGB to the target ("Office
# Theme") value, in clrScheme order (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$targetRGB = @(
    0,           # dk1      000000
    16777215,    # lt1      FFFFFF
    6968388,     # dk2      44546A
    15132391,    # lt2      E7E6E6
    13998939,    # accent1  5B9BD5
    3243501,     # accent2  ED7D31
    10855845,    # accent3  A5A5A5
    49407,       # accent4  FFC000
    12874308,    # accent5  4472C4
    4697456,     # accent6  70AD47
    12673797,    # hlink    0563C1
    7491477      # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $targetRGB[$i - 1]
}
